# Update countries & provincias Spain
# Applies the data refresh captured by the commit diff:
#  - swap the display order of two pairs of tied-rank countries
#    (Santa Lucia / Timor Oriental, and Montserrat / Islas Malvinas)
#  - refresh the "datos actualizados" timestamp
#  - refresh the statistics for a handful of countries whose numbers changed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Timestamp update (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 10:59"

# --- 2. Swap "Santa Lucia" <-> "Timor Oriental" (rows 204/205) ---
# Their statistics are identical, only the country names trade places.
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

# --- 3. Swap "Montserrat" <-> "Islas Malvinas" (rows 214/215) ---
# Here the underlying stats also swap along with the name.
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

# --- 4. Refresh statistics for countries with updated counts ---
# Row 24: Filipinas
$ws.Range("B24").Value = 286743
$ws.Range("C24").Value = 3311
$ws.Range("D24").Value = 229865
$ws.Range("E24").Value = 51894
$ws.Range("G24").Value = 55
$ws.Range("H24").Value = 4984

# Row 26: Indonesia
$ws.Range("B26").Value = 244676
$ws.Range("C26").Value = 3989
$ws.Range("D26").Value = 177327
$ws.Range("E26").Value = 57796
$ws.Range("G26").Value = 105
$ws.Range("H26").Value = 9553

# Row 47: Polonia
$ws.Range("B47").Value = 79240
$ws.Range("C47").Value = 910
$ws.Range("D47").Value = 64302
$ws.Range("E47").Value = 12645
$ws.Range("G47").Value = 11
$ws.Range("H47").Value = 2293

# Row 69: Austria
$ws.Range("B69").Value = 38095
$ws.Range("C69").Value = 621
$ws.Range("D69").Value = 29229
$ws.Range("E69").Value = 8100
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 766

# Row 89: Croacia
$ws.Range("B89").Value = 14922
$ws.Range("C89").Value = 197
$ws.Range("D89").Value = 12536
$ws.Range("E89").Value = 2138
$ws.Range("G89").Value = 4
$ws.Range("H89").Value = 248

# Row 110: Eslovaquia
$ws.Range("B110").Value = 6677
$ws.Range("C110").Value = 131
$ws.Range("D110").Value = 3548
$ws.Range("E110").Value = 3090

# Row 119: Hong Kong
$ws.Range("B119").Value = 5033
$ws.Range("C119").Value = 23
$ws.Range("D119").Value = 4708
$ws.Range("E119").Value = 222

# Row 128: Eslovenia
$ws.Range("B128").Value = 4420
$ws.Range("C128").Value = 111
$ws.Range("D128").Value = 3023
$ws.Range("E128").Value = 1255
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 142

# Row 132: Lituania
$ws.Range("B132").Value = 3744
$ws.Range("C132").Value = 80
$ws.Range("D132").Value = 2198
$ws.Range("E132").Value = 1459

# Row 144: Estonia
$ws.Range("B144").Value = 2924
$ws.Range("C144").Value = 49
$ws.Range("D144").Value = 2377
$ws.Range("E144").Value = 483
